# Aggiornamento dati Serramazzoni al 23 agosto 2021
# Appende le righe 344-357 (dal 2021-08-10 al 2021-08-23) ai dati esistenti.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44418, 1, 5, 58.91363261458702),
    @(44419, 0, 5, 58.91363261458702),
    @(44420, 2, 6, 70.69635913750442),
    @(44421, 3, 9, 106.0445387062566),
    @(44422, 2, 11, 129.6099917520914),
    @(44423, 0, 9, 106.0445387062566),
    @(44424, 1, 9, 106.0445387062566),
    @(44425, 2, 10, 117.827265229174),
    @(44426, 0, 10, 117.827265229174),
    @(44427, 2, 10, 117.827265229174),
    @(44428, 0, 7, 82.47908566042182),
    @(44429, 0, 5, 58.91363261458702),
    @(44430, 7, 12, 141.3927182750088),
    @(44431, 0, 11, 129.6099917520914)
)

$startRow = 344
$formatSource = $ws.Cells.Item($startRow - 1, 1)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# La colonna A (date) usa lo stesso stile (bordo + formato data) delle righe precedenti
$formatSource.Copy()
$destRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($startRow + $data.Count - 1, 1))
$destRange.PasteSpecial(-4122)
